# Applies per-cell text updates to the crypto price/volume table.
# Numeric-looking text values are prefixed with a leading apostrophe so Excel
# keeps storing them as text (matching the original inlineStr string cells)
# instead of silently converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.684.87"
$ws.Range("E2").Value = "  -0.30%  "
$ws.Range("D3").Value = "2.291.89"
$ws.Range("E3").Value = "  -1.05%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'103.81"
$ws.Range("E5").Value = "  +6.92%  "
$ws.Range("D6").Value = "'270.73"
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("D7").Value = "'0.618"
$ws.Range("E7").Value = "  -1.28%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").Value = "'0.608"
$ws.Range("E9").Value = "  -2.00%  "
$ws.Range("D10").Value = "'45.86"
$ws.Range("E10").Value = "  +0.40%  "
$ws.Range("D11").Value = "'0.0934"
$ws.Range("E11").Value = "  -1.11%  "
$ws.Range("D12").Value = "'7.97"
$ws.Range("E12").Value = "  -0.99%  "
$ws.Range("E13").Value = "  +1.23%  "
$ws.Range("D14").Value = "'15.64"
$ws.Range("E14").Value = "  +1.02%  "
$ws.Range("E15").Value = "  -1.42%  "
$ws.Range("D16").Value = "2.299.70"
$ws.Range("E16").Value = "  -0.74%  "
$ws.Range("D17").Value = "43.712.90"
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("E18").Value = "  +0.33%  "
$ws.Range("D19").Value = "'6.26"
$ws.Range("E19").Value = "  -2.42%  "
$ws.Range("D20").Value = "'72.29"
$ws.Range("E20").Value = "  -1.23%  "
$ws.Range("E21").Value = "  +9.93%  "
$ws.Range("D22").Value = "'233.32"
$ws.Range("E22").Value = "  -2.62%  "
$ws.Range("D23").Value = "'2.91"
$ws.Range("E23").Value = "  +14.59%  "
$ws.Range("D24").Value = "'9.18"
$ws.Range("E24").Value = "  -2.92%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").Value = "'11.24"
$ws.Range("E26").Value = "  -1.36%  "
$ws.Range("E27").Value = "  -0.98%  "
$ws.Range("D28").Value = "'39.85"
$ws.Range("E28").Value = "  +4.90%  "
$ws.Range("D29").Value = "'2.23"
$ws.Range("E29").Value = "  -2.26%  "
$ws.Range("D30").Value = "'176.79"
$ws.Range("E30").Value = "  +1.61%  "
$ws.Range("D31").Value = "'21.80"
$ws.Range("E31").Value = "  -2.53%  "
$ws.Range("D32").Value = "'0.0902"
$ws.Range("E32").Value = "  -0.38%  "
$ws.Range("D33").Value = "'5.48"
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("D34").Value = "'4.93"
$ws.Range("E34").Value = "  +11.97%  "
$ws.Range("E35").Value = "  -0.40%  "
$ws.Range("E36").Value = "  +0.28%  "
$ws.Range("E37").Value = "  -2.77%  "
$ws.Range("D38").Value = "'3.55"
$ws.Range("E38").Value = "  +5.63%  "
$ws.Range("E39").Value = "  -4.43%  "
$ws.Range("D40").Value = "'2.33"
$ws.Range("E40").Value = "  -0.96%  "
$ws.Range("D41").Value = "'1.38"
$ws.Range("E41").Value = "  +1.11%  "
$ws.Range("D42").Value = "'12.30"
$ws.Range("E42").Value = "  +0.86%  "
$ws.Range("D43").Value = "'65.50"
$ws.Range("E43").Value = "  +5.07%  "
$ws.Range("B44").Value = "THORChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D44").Value = "'5.25"
$ws.Range("E44").Value = "  -1.59%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "'8.81"
$ws.Range("E45").Value = "  -4.42%  "
$ws.Range("E46").Value = "  -1.34%  "
$ws.Range("E47").Value = "  +2.45%  "
$ws.Range("D48").Value = "'98.86"
$ws.Range("E48").Value = "  -1.62%  "
$ws.Range("D49").Value = "'0.450"
$ws.Range("E49").Value = "  +8.17%  "
$ws.Range("E50").Value = "  +10.99%  "
$ws.Range("D51").Value = "2.515.53"
$ws.Range("E51").Value = "  -1.03%  "
